$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Actual Outcome" (F) and "Status" (G) columns, iteration-1 test results
# Values are written in the specific order below so that the shared-string
# table indices line up with the authored workbook (header->F2->F3->F5->F4
# ->G2->G1->G4->G3->G5 reflects the original authoring order).
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Actual Outcome"
$ws.Range("F2").Value = "The initiative assessment successfully opened upon clicking 'Take Initiative Assessment' in the application"
$ws.Range("F3").Value = "Inside the actual initiative assessment, the user is able to successfully input their answers to each question. The user can successfully submit their assessment by selecting 'submit' and seeing a dialog box open that says 'Your data was submitted'"
$ws.Range("F5").Value = "The app does not send a push notification to the user at the end of the day. However, the end of day assessment was successfully created and can be taken by the user manually through the home page of the app"
$ws.Range("F4").Value = "The app does not send push notifications to the user approximately 5x throughout the day. However, the random assessment was successfully created and can be taken by the user manually through the home page of the app"
$ws.Range("G2").Value = "Success!"
$ws.Range("G1").Value = "Status"
$ws.Range("G4").Value = "Not fully working yet"
$ws.Range("G3").Value = "Success!"
$ws.Range("G5").Value = "Not fully working yet"

# ---------------------------------------------------------------------------
# Row heights (rows grew taller to fit the new wrapped commentary text)
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 120
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 105

# ---------------------------------------------------------------------------
# Column widths - re-fit existing columns B, C, E and size the two new ones
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20.833333333333332
$ws.Columns.Item(3).ColumnWidth = 36.666666666666664
$ws.Columns.Item(5).ColumnWidth = 21.833333333333332
$ws.Columns.Item(6).ColumnWidth = 30.666666666666668
$ws.Columns.Item(7).ColumnWidth = 18.333333333333332

# ---------------------------------------------------------------------------
# Header formatting for the new columns (bold, centered - matches A1:E1)
# ---------------------------------------------------------------------------
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Body formatting - column F (Actual Outcome): centered, wrapped text
# ---------------------------------------------------------------------------
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4160
$ws.Range("F2").WrapText = $true

$ws.Range("F3:F5").HorizontalAlignment = -4108
$ws.Range("F3:F5").VerticalAlignment = -4108
$ws.Range("F3:F5").WrapText = $true

# ---------------------------------------------------------------------------
# Body formatting - column G (Status): centered horizontally & vertically
# ---------------------------------------------------------------------------
$ws.Range("G2:G5").HorizontalAlignment = -4108
$ws.Range("G2:G5").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Existing B2:E2 cells pick up vertical centering to match the rest of the
# table (every other data row already centers vertically)
# ---------------------------------------------------------------------------
$ws.Range("B2:E2").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Selection moves to G3, matching the saved workbook state
# ---------------------------------------------------------------------------
$ws.Range("G3").Select()
